$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report title in B1: drop the period after "6.4.2.1"
$ws.Range("B1").Value = "6.4.2.1 Общий объем забора пресной воды "

# Updated data values (new reporting year figures)
$ws.Range("L5").Value = 8741.9

# L7 used to be a formula (=L5-L8); it is now a plain literal value
$ws.Range("L7").Value = 8483.5

$ws.Range("L14").Value = 1327.6

$ws.Range("L18").Value = 54

# Update the active selection shown when the sheet is opened
$ws.Range("O2").Select() | Out-Null
